# V1.0.5 Done admin photo
#
# Adds new translation rows to the "DATA" sheet of translations.xlsx:
#   - emsTracking                              (row 103)
#   - mailerText/SubjectUpdateOrder (1-3/subj)  (rows 104-107)
#   - mailerText/SubjectFwdOrder    (1-3/subj)  (rows 108-111)
#
# Columns are: A=id  B=en  C=th  D=cn  E=jp  F=es  G=mm  H=la  I=km
# For these new keys, every non-en/th locale simply mirrors the English text.
#
# NOTE: the cell writes below are intentionally ordered so that the
# underlying shared-string table is built up in the same sequence as the
# source workbook (new strings first-used in columns A, then B, then C for
# each new row/group, followed by the duplicate D:I values) rather than in
# simple left-to-right / top-to-bottom order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 103: emsTracking ---------------------------------------------------
$ws.Range("A103").Value = "emsTracking"
$ws.Range("B103").Value = "EMS Tracking"
$ws.Range("C103").Value = "หมายเลขติดตามพัสดุ"

# --- Rows 104-107: "order updated" mailer strings ---------------------------
$ws.Range("A104").Value = "mailerTextUpdateOrder1"
$ws.Range("A105").Value = "mailerTextUpdateOrder2"
$ws.Range("A106").Value = "mailerTextUpdateOrder3"
$ws.Range("A107").Value = "mailerSubjectUpdateOrder"
$ws.Range("B107").Value = "[ICMM2025] Your order has been updated!..."
$ws.Range("B104").Value = "We have updated your order status or detail."
$ws.Range("C107").Value = "[ICMM2025] คำสั่งซื้อของท่านมีการเปลี่ยนสถานะหรือแก้ไข…"
$ws.Range("C104").Value = "เราได้ทำการเปลี่ยนแปลงสถานะหรือรายละเอียดคำสั่งซื้อของท่าน"
$ws.Range("C111").Value = "[ICMM2025] คำสั่งซื้อของท่านมีการเปลี่ยนสถานะ…"

# --- Rows 108-111: "order forwarded to next status" mailer strings ---------
$ws.Range("A108").Value = "mailerTextFwdOrder1"
$ws.Range("A110").Value = "mailerTextFwdOrder3"
$ws.Range("A109").Value = "mailerTextFwdOrder2"
$ws.Range("A111").Value = "mailerSubjectFwdOrder"
$ws.Range("B111").Value = "[ICMM2025] Your order has changed status!..."
$ws.Range("B108").Value = "We have moved your order to the next state."
$ws.Range("C108").Value = "คำสั่งซื้อของท่านถูกทำเดินการไปยังสถานะถัดไป"

# --- Fill remaining locale columns (D:I mirror the English text, B/C reuse
#     existing "check order status" / "thank you" strings where applicable) -
$ws.Range("D103").Value = "EMS Tracking"
$ws.Range("E103").Value = "EMS Tracking"
$ws.Range("F103").Value = "EMS Tracking"
$ws.Range("G103").Value = "EMS Tracking"
$ws.Range("H103").Value = "EMS Tracking"
$ws.Range("I103").Value = "EMS Tracking"

$ws.Range("D104").Value = "We have updated your order status or detail."
$ws.Range("E104").Value = "We have updated your order status or detail."
$ws.Range("F104").Value = "We have updated your order status or detail."
$ws.Range("G104").Value = "We have updated your order status or detail."
$ws.Range("H104").Value = "We have updated your order status or detail."
$ws.Range("I104").Value = "We have updated your order status or detail."

$ws.Range("B105").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("C105").Value = 'ท่านสามารถตรวจสอบสถานะคำสั่งซื้อได้โดยคลิกที่ "ตรวจสอบสถานะคำสั่งซื้อ" และกรอกหมายเลขคำสั่งซื้อของท่านคือ :'
$ws.Range("D105").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("E105").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("F105").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("G105").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("H105").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("I105").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'

$ws.Range("B106").Value = "Thank you for shopping with us."
$ws.Range("C106").Value = "ขอขอบพระคุณที่เลือกซื้อสินค้ากับเรา"
$ws.Range("D106").Value = "Thank you for shopping with us."
$ws.Range("E106").Value = "Thank you for shopping with us."
$ws.Range("F106").Value = "Thank you for shopping with us."
$ws.Range("G106").Value = "Thank you for shopping with us."
$ws.Range("H106").Value = "Thank you for shopping with us."
$ws.Range("I106").Value = "Thank you for shopping with us."

$ws.Range("D107").Value = "[ICMM2025] Your order has been updated!..."
$ws.Range("E107").Value = "[ICMM2025] Your order has been updated!..."
$ws.Range("F107").Value = "[ICMM2025] Your order has been updated!..."
$ws.Range("G107").Value = "[ICMM2025] Your order has been updated!..."
$ws.Range("H107").Value = "[ICMM2025] Your order has been updated!..."
$ws.Range("I107").Value = "[ICMM2025] Your order has been updated!..."

$ws.Range("D108").Value = "We have moved your order to the next state."
$ws.Range("E108").Value = "We have moved your order to the next state."
$ws.Range("F108").Value = "We have moved your order to the next state."
$ws.Range("G108").Value = "We have moved your order to the next state."
$ws.Range("H108").Value = "We have moved your order to the next state."
$ws.Range("I108").Value = "We have moved your order to the next state."

$ws.Range("B109").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("C109").Value = 'ท่านสามารถตรวจสอบสถานะคำสั่งซื้อได้โดยคลิกที่ "ตรวจสอบสถานะคำสั่งซื้อ" และกรอกหมายเลขคำสั่งซื้อของท่านคือ :'
$ws.Range("D109").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("E109").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("F109").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("G109").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("H109").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'
$ws.Range("I109").Value = 'You can check your order status by clicking "Check Order Status" and entering your order no. :'

$ws.Range("B110").Value = "Thank you for shopping with us."
$ws.Range("C110").Value = "ขอขอบพระคุณที่เลือกซื้อสินค้ากับเรา"
$ws.Range("D110").Value = "Thank you for shopping with us."
$ws.Range("E110").Value = "Thank you for shopping with us."
$ws.Range("F110").Value = "Thank you for shopping with us."
$ws.Range("G110").Value = "Thank you for shopping with us."
$ws.Range("H110").Value = "Thank you for shopping with us."
$ws.Range("I110").Value = "Thank you for shopping with us."

$ws.Range("D111").Value = "[ICMM2025] Your order has changed status!..."
$ws.Range("E111").Value = "[ICMM2025] Your order has changed status!..."
$ws.Range("F111").Value = "[ICMM2025] Your order has changed status!..."
$ws.Range("G111").Value = "[ICMM2025] Your order has changed status!..."
$ws.Range("H111").Value = "[ICMM2025] Your order has changed status!..."
$ws.Range("I111").Value = "[ICMM2025] Your order has changed status!..."

# --- Column C is now much wider (long Thai mailer strings) ----------------
$ws.Columns.Item(3).ColumnWidth = 47

# --- Update the view: scroll down and select the newly-added rows ---------
$ws.Range("D104:I111").Select()
$excel.ActiveWindow.ScrollRow = 93
$excel.ActiveWindow.ScrollColumn = 1
